# Applies the crypto price/volume table update described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.329.61"
$ws.Range("E2").Value = "  -1.36%  "

# Row 3
$ws.Range("D3").Value = "2.520.79"
$ws.Range("E3").Value = "  -0.94%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.73"
$ws.Range("E5").Value = "  +2.67%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.73"
$ws.Range("E6").Value = "  -5.67%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -1.08%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -3.59%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.60"
$ws.Range("E10").Value = "  -4.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  -2.49%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.56"
$ws.Range("E12").Value = "  -2.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.113"
$ws.Range("E13").Value = "  -0.47%  "

# Row 14
$ws.Range("D14").Value = "2.909.80"
$ws.Range("E14").Value = "  -0.82%  "

# Row 15
$ws.Range("D15").Value = "2.521.05"
$ws.Range("E15").Value = "  -0.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.36"
$ws.Range("E16").Value = "  +0.64%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("E17").Value = "  -3.63%  "

# Row 18
$ws.Range("D18").Value = "42.435.04"
$ws.Range("E18").Value = "  -1.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  -2.44%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.53"
$ws.Range("E20").Value = "  -0.58%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0955"
$ws.Range("E21").Value = "  -3.63%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.43"
$ws.Range("E22").Value = "  -1.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.99"
$ws.Range("E23").Value = "  -2.21%  "

# Row 24
$ws.Range("E24").Value = "  -0.84%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").Value = "  -3.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.37"
$ws.Range("E26").Value = "  -5.03%  "

# Row 27
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("E28").Value = "  +3.22%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.75"
$ws.Range("E29").Value = "  -0.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.09"
$ws.Range("E30").Value = "  -1.25%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.87"
$ws.Range("E31").Value = "  -5.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.23"
$ws.Range("E32").Value = "  -0.99%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.28"
$ws.Range("E33").Value = "  +2.38%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.11"
$ws.Range("E34").Value = "  -1.64%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.28"
$ws.Range("E35").Value = "  -0.66%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0778"
$ws.Range("E36").Value = "  -2.97%  "

# Row 37
$ws.Range("E37").Value = "  -1.03%  "

# Row 38
$ws.Range("E38").Value = "  -4.88%  "

# Row 39
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.84"
$ws.Range("E39").Value = "  -1.50%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.118"
$ws.Range("E40").Value = "  -1.81%  "

# Row 41
$ws.Range("E41").Value = "  +10.78%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.34%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.78"
$ws.Range("E43").Value = "  -3.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("E44").Value = "  -5.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0298"
$ws.Range("E45").Value = "  -2.27%  "

# Row 46
$ws.Range("D46").Value = "2.013.05"
$ws.Range("E46").Value = "  -2.91%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "83.94"
$ws.Range("E47").Value = "  -2.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.78"
$ws.Range("E48").Value = "  -2.88%  "

# Row 49
$ws.Range("D49").Value = "2.767.59"
$ws.Range("E49").Value = "  -0.85%  "

# Row 50
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.37"
$ws.Range("E50").Value = "  -1.79%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.57"
$ws.Range("E51").Value = "  -1.82%  "
